$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - copy formatting from the existing header cell G1
# (bold, centered, bordered style) then set its text to "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column (H2:H9)
$values = @(1, 0, 0, 0, 0, 1, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
